# Updates to lemmatizer for more accuracy
# - Refreshes the lemmatized "Terms" text for each topic row (D2:D6)
# - Refreshes the recomputed Score for Negativity / Percent Dominance values (A2:B6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated term lists (Terms column, D2:D6)
$ws.Range("D2").Value = "call,phone,time,tell,day,try,say,harass,go,number,work,make,know,even,give,back,help,ask,people,speak"
$ws.Range("D3").Value = "interest,year,pay,high,principal,amount,rate,payment,leave,balance,make,go,consolidate,charge,total,borrow,money,take,apply,month"
$ws.Range("D4").Value = "payment,late,fee,account,check,send,bank,charge,month,statement,receive,make,apply,letter,mail,cash,due,copy,never,amount"
$ws.Range("D5").Value = "report,credit,reply,list,however,remove,year,ask,debt,another,since,pay,write,right,score,keep,agency,lose,greedy,could"
$ws.Range("D6").Value = "present,already,refuse,bill,pay,original,balance,threat,send,legal,lower,college,son,want,talk,help,correspondence,able,cosigner,clue"

# Updated Score for Negativity (A2:A6) and Percent Dominance (B2:B6)
$ws.Range("A2").Value = 0.3192080265513397
$ws.Range("B2").Value = 0.5274809160305344

$ws.Range("A3").Value = -0.294468947278769
$ws.Range("B3").Value = 0.1740458015267176

$ws.Range("A4").Value = -0.6584019662348324
$ws.Range("B4").Value = 0.2198473282442748

$ws.Range("A5").Value = -0.3119978868290287
$ws.Range("B5").Value = 0.06106870229007633

$ws.Range("A6").Value = -0.5005420399354448
$ws.Range("B6").Value = 0.01755725190839695
